$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 7144361.5
$ws.Range("I6").Value = 1755.1666
$ws.Range("J6").Value = 50000000
$ws.Range("K6").Value = 5265.4998
$ws.Range("L6").Value = 150000000
$ws.Range("M6").Value = -5153.4998
$ws.Range("N6").Value = -150000224

$ws.Range("H33").Value = 7679.0713
$ws.Range("I33").Value = 72.57143000000001
$ws.Range("J33").Value = 15285.571
$ws.Range("K33").Value = 72.57143000000001
$ws.Range("L33").Value = 15285.571
$ws.Range("M33").Value = 156.42857
$ws.Range("N33").Value = -15743.571

$ws.Range("H100").Value = 10418491
$ws.Range("I100").Value = 12821747
$ws.Range("J100").Value = 4381
$ws.Range("K100").Value = 12821747
$ws.Range("L100").Value = 4381
$ws.Range("M100").Value = -12821206
$ws.Range("N100").Value = -5463

$ws.Range("H137").Value = 923.6667
$ws.Range("I137").Value = 811.4286
$ws.Range("J137").Value = 1021.875
$ws.Range("K137").Value = 2434.2858
$ws.Range("L137").Value = 3065.625
$ws.Range("M137").Value = 115.7142000000003
$ws.Range("N137").Value = -8165.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1704.5769
$ws.Range("I2").Value = 746.64703
$ws.Range("J2").Value = 3514
$ws.Range("K2").Value = 746.64703
$ws.Range("L2").Value = 3514
$ws.Range("M2").Value = -633.64703
$ws.Range("N2").Value = -3740

$ws.Range("H74").Value = 685.1786
$ws.Range("I74").Value = 624.85
$ws.Range("J74").Value = 836
$ws.Range("K74").Value = 624.85
$ws.Range("L74").Value = 836
$ws.Range("M74").Value = 249.15
$ws.Range("N74").Value = -2584

$ws.Range("H77").Value = 685.1786
$ws.Range("I77").Value = 624.85
$ws.Range("J77").Value = 836
$ws.Range("K77").Value = 3124.25
$ws.Range("L77").Value = 4180
$ws.Range("M77").Value = 1243.75
$ws.Range("N77").Value = -12916

$ws.Range("H116").Value = 1704.5769
$ws.Range("I116").Value = 746.64703
$ws.Range("J116").Value = 3514
$ws.Range("K116").Value = 746.64703
$ws.Range("L116").Value = 3514
$ws.Range("M116").Value = 1547.35297
$ws.Range("N116").Value = -8102

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1704.5769
$ws.Range("I3").Value = 746.64703
$ws.Range("J3").Value = 3514
$ws.Range("K3").Value = 746.64703
$ws.Range("L3").Value = 3514
$ws.Range("M3").Value = -632.64703
$ws.Range("N3").Value = -3742

$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21372

$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66864

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2107.347
$ws.Range("I31").Value = 1609.2307
$ws.Range("K31").Value = 1609.2307
$ws.Range("M31").Value = -1314.2307

$ws.Range("H34").Value = 2107.347
$ws.Range("I34").Value = 1609.2307
$ws.Range("K34").Value = 1609.2307
$ws.Range("M34").Value = -1407.2307

$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1706

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1077.721
$ws.Range("I113").Value = 953.4375
$ws.Range("J113").Value = 1151.3704
$ws.Range("K113").Value = 2860.3125
$ws.Range("L113").Value = 3454.1112
$ws.Range("M113").Value = -690.3125
$ws.Range("N113").Value = -7794.1112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 3300
$ws.Range("J40").Value = 3300
$ws.Range("L40").Value = 3300
$ws.Range("N40").Value = -3602

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1895.3043
$ws.Range("I40").Value = 1866.2222
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1866.2222
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1730.2222
$ws.Range("N40").Value = -2272

$ws.Range("H68").Value = 2028.7715
$ws.Range("J68").Value = 2300.5386
$ws.Range("L68").Value = 2300.5386
$ws.Range("N68").Value = -3798.5386

$ws.Range("H71").Value = 2028.7715
$ws.Range("J71").Value = 2300.5386
$ws.Range("L71").Value = 11502.693
$ws.Range("N71").Value = -18990.693

$ws.Range("H82").Value = 2652.1428
$ws.Range("I82").Value = 3016.25
$ws.Range("J82").Value = 2166.6667
$ws.Range("K82").Value = 3016.25
$ws.Range("L82").Value = 2166.6667
$ws.Range("M82").Value = -2655.25
$ws.Range("N82").Value = -2888.6667

$ws.Range("H85").Value = 2652.1428
$ws.Range("I85").Value = 3016.25
$ws.Range("J85").Value = 2166.6667
$ws.Range("K85").Value = 3016.25
$ws.Range("L85").Value = 2166.6667
$ws.Range("M85").Value = -1768.25
$ws.Range("N85").Value = -4662.6667

$ws.Range("H93").Value = 1305.9736
$ws.Range("I93").Value = 1155.84
$ws.Range("J93").Value = 1594.6923
$ws.Range("K93").Value = 1155.84
$ws.Range("L93").Value = 1594.6923
$ws.Range("M93").Value = 92.16000000000008
$ws.Range("N93").Value = -4090.6923

$ws.Range("H134").Value = 48707.69
$ws.Range("J134").Value = 48707.69
$ws.Range("L134").Value = 48707.69
$ws.Range("N134").Value = -58847.69

$ws.Range("H136").Value = 4481
$ws.Range("I136").Value = 4150.697
$ws.Range("J136").Value = 5571
$ws.Range("K136").Value = 12452.091
$ws.Range("L136").Value = 16713
$ws.Range("M136").Value = -9902.091
$ws.Range("N136").Value = -21813

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 4817.8184
$ws.Range("J49").Value = 4999.6
$ws.Range("L49").Value = 4999.6
$ws.Range("N49").Value = -5459.6

$ws.Range("H52").Value = 10000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
